# Update the "Short Name" column (column B) on the PolicyLevers sheet:
# insert a space before every hyphen so that e.g.
#   "Agriculture, Land Use, and Water- Desalination Energy Efficiency Standards"
# becomes
#   "Agriculture, Land Use, and Water - Desalination Energy Efficiency Standards"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PolicyLevers")

# Data (excluding the header in B1) runs from B2 down to the last used row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
$firstCell = $ws.Cells.Item(2, 2)
$lastCell = $ws.Cells.Item($lastRow, 2)
$rng = $ws.Range($firstCell, $lastCell)
[void]$rng.Replace("-", " -", 2, 1, $false, $false, $false)

# Leave the selection where the user last left it after editing the column.
$ws.Activate()
[void]$ws.Range("B18").Select()
